# Add the team's season record (Wins / Losses / Ties) as three new
# columns (AD, AE, AF) appended after the existing data (which ends at
# column AC), for every row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells should look like the rest of the header row (bold,
# bordered, centered) - copy the formatting from the last existing header
# cell (AC1) onto the three new header cells, then set their captions.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2 through 42) gets the same team season record.
$ws.Range("AD2:AD42").Value = 96
$ws.Range("AE2:AE42").Value = 66
$ws.Range("AF2:AF42").Value = 0
